$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subjects")

$ws.Range("D3").Value = 15
$ws.Range("D4").Value = 12
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 9
$ws.Range("D8").Value = 12
$ws.Range("D11").Value = 12
$ws.Range("D12").Value = 12
$ws.Range("D15").Value = 27
$ws.Range("D16").Value = 27
$ws.Range("D18").Value = 21

$ws.Activate()
$ws.Range("E20").Select()
